$wb = $excel.ActiveWorkbook
$wsChild = $wb.Worksheets.Item("Child")

# Update the computed "address" (time-travel coordinate) values on the Child sheet
$wsChild.Range("D2").Value = "-4.17,-3.8"
$wsChild.Range("D3").Value = "1.95,-8.2"
$wsChild.Range("D4").Value = "-7.38,-6.34"
$wsChild.Range("D5").Value = "-2.83,7.67"
$wsChild.Range("D6").Value = "2.15,-5.14"
$wsChild.Range("D7").Value = "7.41,-6.16"
$wsChild.Range("D8").Value = "-6.44,3.18"
$wsChild.Range("D9").Value = "-2.87,-2.03"
$wsChild.Range("D10").Value = "-2.69,6.26"
$wsChild.Range("D11").Value = "-9.1,7.31"
$wsChild.Range("D12").Value = "1.4,-5.2"
$wsChild.Range("D13").Value = "9.68,-3.45"
$wsChild.Range("D14").Value = "3.22,4.01"
$wsChild.Range("D15").Value = "-3.16,-7.95"
$wsChild.Range("D16").Value = "-0.99,-5.37"
$wsChild.Range("D17").Value = "-9.16,-3.53"
$wsChild.Range("D18").Value = "4.71,-5.7"
$wsChild.Range("D19").Value = "-3.08,-1.04"
$wsChild.Range("D20").Value = "6.6,-1.08"
$wsChild.Range("D21").Value = "-1.93,9.03"
$wsChild.Range("D22").Value = "3.24,-0.78"

# Widen column E (address) so the longer values are fully visible
$wsChild.Columns.Item(5).ColumnWidth = 29.83

# Make the Child sheet the active tab/sheet and select E12 there
# (moves focus away from the School sheet, which loses tabSelected)
$wsChild.Activate()
$wsChild.Range("E12").Select()
